$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/55c93a3534181f48ae13d92fbc93946ec87b392f/e2e/508b30cf-0b27-442e-b887-a496708c3359.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f11cdf907024b6b4c22349eb02f9a494e13d2b/e2e/508b30cf-0b27-442e-b887-a496708c3359.md."
$targetMdAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01f11cdf907024b6b4c22349eb02f9a494e13d2b/e2e/508b30cf-0b27-442e-b887-a496708c3359.md"
$targetMdDisplay = "508b30cf-0b27-442e-b887-a496708c3359.md"

# ----- zh-cn sheet -----
$ws = $wb.Worksheets.Item("zh-cn")

# widen the Error Detail column (col 16 / P) to match the other "wide" columns (e.g. col A)
$wideWidth = $ws.Columns.Item(1).ColumnWidth
$ws.Columns.Item(16).ColumnWidth = $wideWidth

# Latest Target File (I5) now links to the handed-back markdown file
$ws.Hyperlinks.Add($ws.Range("I5"), $targetMdAddress, "", "", $targetMdDisplay)
$ws.Range("I5").Font.Underline = 2
$ws.Range("I5").Font.Color = 15570276

# Latest Handback File (J5)
$ws.Range("J5").Value = "508b30cf-0b27-442e-b887-a496708c3359.a51854df91fa07fae867c18bfb2953576e7f11c4.zh-cn.xlf"

# Latest Handback DateTime (K5)
$ws.Range("K5").Value = "2016-10-20 09:03:06"

# Error Detail (P5)
$ws.Range("P5").Value = $errorDetail

# ----- de-de sheet -----
$ws2 = $wb.Worksheets.Item("de-de")

$wideWidth2 = $ws2.Columns.Item(1).ColumnWidth
$ws2.Columns.Item(16).ColumnWidth = $wideWidth2

$ws2.Hyperlinks.Add($ws2.Range("I5"), $targetMdAddress, "", "", $targetMdDisplay)
$ws2.Range("I5").Font.Underline = 2
$ws2.Range("I5").Font.Color = 15570276

$ws2.Range("J5").Value = "508b30cf-0b27-442e-b887-a496708c3359.a51854df91fa07fae867c18bfb2953576e7f11c4.de-de.xlf"

$ws2.Range("K5").Value = "2016-10-20 09:03:23"

$ws2.Range("P5").Value = $errorDetail

$wb.Save()
